$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in "yes" for the "Has this functionality been implemented?" column
# for every requirement row (2-21).
for ($r = 2; $r -le 21; $r++) {
    $ws.Cells.Item($r, 2).Value = "yes"
}

# New blank row 22 that only carries the wrap-text style in column C.
$ws.Range("C22").WrapText = $true

# Highlight/re-font the "view and edit driver status" requirement (row 7).
$ws.Range("D7").Font.Name = "新細明體 (本文)"

# Widen columns C & D to fit the full requirement text, and shrink the
# header row a touch.
$ws.Columns.Item(3).ColumnWidth = 25.330729166666668
$ws.Columns.Item(4).ColumnWidth = 88.49869791666667
$ws.Rows.Item(1).RowHeight = 48

# Zoom in a bit and move the active selection, matching the saved view state.
$ws.Application.ActiveWindow.Zoom = 125
[void]$ws.Range("C20").Select()
